$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 7's A:H block entirely (no longer present after the rerun)
$ws.Range("A7:H7").Clear()

# Update anchor-word cells and their stats for the rerun with a larger dataset
$ws.Range("A1").Value = "negative"
$ws.Range("J1").Value = "positive"
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "crude"
$ws.Range("B3").Value = 0.7941176470588235
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 27
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $false
$ws.Range("H3").Value = 7
$ws.Range("J3").Value = "best"
$ws.Range("K3").Value = 0.9661016949152542
$ws.Range("L3").Value = 57
$ws.Range("M3").Value = 57
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $false
$ws.Range("Q3").Value = 2
$ws.Range("A4").Value = "crisis"
$ws.Range("B4").Value = 0.6232876712328768
$ws.Range("C4").Value = 182
$ws.Range("D4").Value = 182
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = 110
$ws.Range("J4").Value = "happy"
$ws.Range("K4").Value = 0.9615384615384616
$ws.Range("L4").Value = 25
$ws.Range("M4").Value = 25
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $false
$ws.Range("Q4").Value = 1
$ws.Range("A5").Value = "panic"
$ws.Range("B5").Value = 0.2151162790697674
$ws.Range("C5").Value = 111
$ws.Range("D5").Value = 111
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 405
$ws.Range("J5").Value = "love"
$ws.Range("K5").Value = 0.8913043478260869
$ws.Range("L5").Value = 41
$ws.Range("M5").Value = 41
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $false
$ws.Range("Q5").Value = 5
$ws.Range("A6").Value = "sc"
$ws.Range("B6").Value = 0.164021164021164
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 158
$ws.Range("J6").Value = "interesting"
$ws.Range("K6").Value = 0.8787878787878788
$ws.Range("L6").Value = 29
$ws.Range("M6").Value = 29
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $false
$ws.Range("Q6").Value = 4
$ws.Range("J7").Value = "great"
$ws.Range("K7").Value = 0.8482142857142857
$ws.Range("L7").Value = 95
$ws.Range("M7").Value = 95
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = 17
$ws.Range("J8").Value = "thanks"
$ws.Range("K8").Value = 0.8292682926829268
$ws.Range("L8").Value = 68
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $false
$ws.Range("Q8").Value = 14
$ws.Range("J9").Value = "positive"
$ws.Range("K9").Value = 0.8103448275862069
$ws.Range("L9").Value = 47
$ws.Range("M9").Value = 47
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $false
$ws.Range("Q9").Value = 11
$ws.Range("J10").Value = "thank"
$ws.Range("K10").Value = 0.78125
$ws.Range("L10").Value = 100
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $false
$ws.Range("Q10").Value = 28
$ws.Range("J11").Value = "confidence"
$ws.Range("K11").Value = 0.7777777777777778
$ws.Range("L11").Value = 28
$ws.Range("M11").Value = 28
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $false
$ws.Range("Q11").Value = 8
$ws.Range("J12").Value = "special"
$ws.Range("K12").Value = 0.7777777777777778
$ws.Range("L12").Value = 28
$ws.Range("M12").Value = 28
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $false
$ws.Range("Q12").Value = 8
$ws.Range("J13").Value = "free"
$ws.Range("K13").Value = 0.7666666666666667
$ws.Range("L13").Value = 92
$ws.Range("M13").Value = 92
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $false
$ws.Range("Q13").Value = 28
$ws.Range("J14").Value = "safe"
$ws.Range("K14").Value = 0.7323943661971831
$ws.Range("L14").Value = 104
$ws.Range("M14").Value = 104
$ws.Range("N14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = $false
$ws.Range("Q14").Value = 38
$ws.Range("J15").Value = "good"
$ws.Range("K15").Value = 0.70625
$ws.Range("L15").Value = 113
$ws.Range("M15").Value = 113
$ws.Range("N15").Value = 1
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = $false
$ws.Range("Q15").Value = 47
$ws.Range("J16").Value = "support"
$ws.Range("K16").Value = 0.6981132075471698
$ws.Range("L16").Value = 74
$ws.Range("M16").Value = 74
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = $false
$ws.Range("Q16").Value = 32
$ws.Range("J17").Value = "better"
$ws.Range("K17").Value = 0.6666666666666666
$ws.Range("L17").Value = 42
$ws.Range("M17").Value = 42
$ws.Range("N17").Value = 1
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = $false
$ws.Range("Q17").Value = 21
$ws.Range("J18").Value = "heroes"
$ws.Range("K18").Value = 0.6595744680851063
$ws.Range("L18").Value = 31
$ws.Range("M18").Value = 31
$ws.Range("N18").Value = 1
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = $false
$ws.Range("Q18").Value = 16
$ws.Range("J19").Value = "safety"
$ws.Range("K19").Value = 0.6470588235294118
$ws.Range("L19").Value = 33
$ws.Range("M19").Value = 33
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = $false
$ws.Range("Q19").Value = 18
$ws.Range("J20").Value = "fresh"
$ws.Range("K20").Value = 0.6458333333333334
$ws.Range("L20").Value = 31
$ws.Range("M20").Value = 31
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = $false
$ws.Range("Q20").Value = 17
$ws.Range("J21").Value = "well"
$ws.Range("K21").Value = 0.6382978723404256
$ws.Range("L21").Value = 60
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = $false
$ws.Range("Q21").Value = 34
$ws.Range("J22").Value = "hand"
$ws.Range("K22").Value = 0.5561357702349869
$ws.Range("L22").Value = 213
$ws.Range("M22").Value = 213
$ws.Range("N22").Value = 1
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = $false
$ws.Range("Q22").Value = 170
$ws.Range("J23").Value = "relief"
$ws.Range("K23").Value = 0.52
$ws.Range("L23").Value = 26
$ws.Range("M23").Value = 26
$ws.Range("N23").Value = 1
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = $false
$ws.Range("Q23").Value = 24
$ws.Range("J24").Value = "care"
$ws.Range("K24").Value = 0.4719101123595505
$ws.Range("L24").Value = 42
$ws.Range("M24").Value = 42
$ws.Range("N24").Value = 1
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = $false
$ws.Range("Q24").Value = 47
$ws.Range("J25").Value = "like"
$ws.Range("K25").Value = 0.4647058823529412
$ws.Range("L25").Value = 158
$ws.Range("M25").Value = 158
$ws.Range("N25").Value = 1
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = $false
$ws.Range("Q25").Value = 182
$ws.Range("J26").Value = "help"
$ws.Range("K26").Value = 0.423728813559322
$ws.Range("L26").Value = 125
$ws.Range("M26").Value = 125
$ws.Range("N26").Value = 1
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = $false
$ws.Range("Q26").Value = 170
$ws.Range("J27").Value = "please"
$ws.Range("K27").Value = 0.3598326359832636
$ws.Range("L27").Value = 86
$ws.Range("M27").Value = 86
$ws.Range("N27").Value = 1
$ws.Range("O27").Value = 0
$ws.Range("P27").Value = $false
$ws.Range("Q27").Value = 153
$ws.Range("J28").Value = "co"
$ws.Range("K28").Value = 0.009980682549903412
$ws.Range("L28").Value = 31
$ws.Range("M28").Value = 32
$ws.Range("N28").Value = 0.97
$ws.Range("O28").Value = 0.03000000000000003
$ws.Range("P28").Value = $true
$ws.Range("Q28").Value = 3075
